$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Temporarily force text format on the Price/Volume columns so that
# values such as "36.891.32" or "244.99" are stored as text, matching
# the original inline-string cells, instead of being auto-parsed as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "36.891.32"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.051.19"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D5").Value = "244.99"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "0.654"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D7").Value = "57.30"
$ws.Range("E7").Value = "  -2.91%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "59.03"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("E10").Value = "  -3.88%  "
$ws.Range("D11").Value = "0.0775"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "15.07"
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("D14").Value = "0.874"
$ws.Range("E14").Value = "  +4.83%  "
$ws.Range("D15").Value = "2.351.24"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "5.56"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").Value = "2.062.61"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "36.877.94"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "17.40"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").Value = "73.05"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").Value = "235.58"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "10.01"
$ws.Range("E26").Value = "  +6.85%  "
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").Value = "168.76"
$ws.Range("D29").Value = "20.08"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "5.44"
$ws.Range("E30").Value = "  +14.32%  "
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").Value = "4.80"
$ws.Range("E33").Value = "  +6.24%  "
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("E35").Value = "  +5.65%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +5.63%  "
$ws.Range("D38").Value = "0.0853"
$ws.Range("E38").Value = "  -4.90%  "
$ws.Range("D39").Value = "1.30"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  -6.70%  "
$ws.Range("D42").Value = "4.84"
$ws.Range("E42").Value = "  -5.63%  "
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "0.0949"
$ws.Range("E44").Value = "  -10.85%  "
$ws.Range("D45").Value = "96.68"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "16.69"
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("D47").Value = "1.303.18"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "2.237.52"
$ws.Range("E51").Value = "  +0.08%  "

# Restore the original (default) cell formatting.
$dataRange.ClearFormats()
